$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = "'" + '29.374.95'
$ws.Cells.Item(2, 4).Style = "Normal"
$ws.Cells.Item(2, 5).Value = '  -0.12%  '

# Row 3
$ws.Cells.Item(3, 4).Value = "'" + '1.844.48'
$ws.Cells.Item(3, 4).Style = "Normal"
$ws.Cells.Item(3, 5).Value = '  -0.29%  '

# Row 4
$ws.Cells.Item(4, 4).Value = "'" + '0.9976'
$ws.Cells.Item(4, 4).Style = "Normal"
$ws.Cells.Item(4, 5).Value = '  -0.35%  '

# Row 5
$ws.Cells.Item(5, 4).Value = "'" + '240.54'
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = '  +0.04%  '

# Row 6
$ws.Cells.Item(6, 4).Value = "'" + '0.6322'
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = '  +0.65%  '

# Row 7
$ws.Cells.Item(7, 4).Value = "'" + '0.9990'
$ws.Cells.Item(7, 4).Style = "Normal"
$ws.Cells.Item(7, 5).Value = '  -0.27%  '

# Row 8
$ws.Cells.Item(8, 4).Value = "'" + '0.07502'
$ws.Cells.Item(8, 4).Style = "Normal"
$ws.Cells.Item(8, 5).Value = '  -1.72%  '

# Row 9
$ws.Cells.Item(9, 5).Value = '  -0.01%  '

# Row 10
$ws.Cells.Item(10, 4).Value = "'" + '24.44'
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Cells.Item(10, 5).Value = '  -1.16%  '

# Row 11
$ws.Cells.Item(11, 4).Value = "'" + '0.07725'
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Cells.Item(11, 5).Value = '  -0.32%  '

# Row 12
$ws.Cells.Item(12, 4).Value = "'" + '1.845.01'
$ws.Cells.Item(12, 4).Style = "Normal"
$ws.Cells.Item(12, 5).Value = '  -2.29%  '

# Row 13
$ws.Cells.Item(13, 4).Value = "'" + '5.006'
$ws.Cells.Item(13, 4).Style = "Normal"
$ws.Cells.Item(13, 5).Value = '  -0.58%  '

# Row 14
$ws.Cells.Item(14, 4).Value = "'" + '0.6793'
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Cells.Item(14, 5).Value = '  +0.05%  '

# Row 15
$ws.Cells.Item(15, 4).Value = "'" + '0.00001038'
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Cells.Item(15, 5).Value = '  -2.87%  '

# Row 16
$ws.Cells.Item(16, 5).Value = '  -1.35%  '

# Row 17
$ws.Cells.Item(17, 4).Value = "'" + '2.104.40'
$ws.Cells.Item(17, 4).Style = "Normal"
$ws.Cells.Item(17, 5).Value = '  -3.77%  '

# Row 18
$ws.Cells.Item(18, 4).Value = "'" + '6.151'
$ws.Cells.Item(18, 4).Style = "Normal"
$ws.Cells.Item(18, 5).Value = '  -0.21%  '

# Row 19
$ws.Cells.Item(19, 4).Value = "'" + '29.384.83'
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Cells.Item(19, 5).Value = '  -0.23%  '

# Row 20
$ws.Cells.Item(20, 4).Value = "'" + '228.94'
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(20, 5).Value = '  +1.12%  '

# Row 21
$ws.Cells.Item(21, 4).Value = "'" + '12.33'
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Cells.Item(21, 5).Value = '  +0.05%  '

# Row 22
$ws.Cells.Item(22, 4).Value = "'" + '0.9990'
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).Value = '  -0.24%  '

# Row 23
$ws.Cells.Item(23, 4).Value = "'" + '7.443'
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Cells.Item(23, 5).Value = '  -0.26%  '

# Row 24
$ws.Cells.Item(24, 4).Value = "'" + '0.9988'
$ws.Cells.Item(24, 4).Style = "Normal"

# Row 25
$ws.Cells.Item(25, 4).Value = "'" + '158.90'
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).Value = '  +0.72%  '

# Row 26
$ws.Cells.Item(26, 5).Value = '  -0.14%  '

# Row 27
$ws.Cells.Item(27, 4).Value = "'" + '8.420'
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Cells.Item(27, 5).Value = '  +0.02%  '

# Row 28
$ws.Cells.Item(28, 5).Value = '  -0.65%  '

# Row 29
$ws.Cells.Item(29, 4).Value = "'" + '0.06433'
$ws.Cells.Item(29, 4).Style = "Normal"
$ws.Cells.Item(29, 5).Value = '  +15.07%  '

# Row 30
$ws.Cells.Item(30, 4).Value = "'" + '1.384'
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Cells.Item(30, 5).Value = '  +0.08%  '

# Row 31
$ws.Cells.Item(31, 5).Value = '  +0.65%  '

# Row 32
$ws.Cells.Item(32, 4).Value = "'" + '4.094'
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Cells.Item(32, 5).Value = '  -0.78%  '

# Row 33
$ws.Cells.Item(33, 4).Value = "'" + '4.061'
$ws.Cells.Item(33, 4).Style = "Normal"
$ws.Cells.Item(33, 5).Value = '  +0.00%  '

# Row 34
$ws.Cells.Item(34, 5).Value = '  -1.12%  '

# Row 35
$ws.Cells.Item(35, 5).Value = '  -1.84%  '

# Row 36
$ws.Cells.Item(36, 4).Value = "'" + '0.7007'
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Cells.Item(36, 5).Value = '  +0.70%  '

# Row 37
$ws.Cells.Item(37, 4).Value = "'" + '2.577'
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Cells.Item(37, 5).Value = '  -0.44%  '

# Row 38
$ws.Cells.Item(38, 4).Value = "'" + '2.835'
$ws.Cells.Item(38, 4).Style = "Normal"
$ws.Cells.Item(38, 5).Value = '  +4.02%  '

# Row 39
$ws.Cells.Item(39, 4).Value = "'" + '1.255.72'
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Cells.Item(39, 5).Value = '  +2.19%  '

# Row 41
$ws.Cells.Item(41, 4).Value = "'" + '6.597'
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Cells.Item(41, 5).Value = '  +2.91%  '

# Row 42
$ws.Cells.Item(42, 4).Value = "'" + '0.9060'
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(42, 5).Value = '  +0.04%  '

# Row 43
$ws.Cells.Item(43, 4).Value = "'" + '0.9985'
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Cells.Item(43, 5).Value = '  -0.31%  '

# Row 44
$ws.Cells.Item(44, 4).Value = "'" + '2.007.76'
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(44, 5).Value = '  -18.41%  '

# Row 45
$ws.Cells.Item(45, 4).Value = "'" + '101.35'
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(45, 5).Value = '  -0.36%  '

# Row 46
$ws.Cells.Item(46, 4).Value = "'" + '66.36'
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(46, 5).Value = '  +0.64%  '

# Row 47
$ws.Cells.Item(47, 2).Value = 'Algorand'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Cells.Item(47, 4).Value = "'" + '0.1178'
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 5).Value = '  +2.86%  '

# Row 48
$ws.Cells.Item(48, 2).Value = 'Aptos'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Cells.Item(48, 4).Value = "'" + '7.051'
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Cells.Item(48, 5).Value = '  -1.72%  '

# Row 49
$ws.Cells.Item(49, 2).Value = 'BabyDogeCoin'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Cells.Item(49, 4).Value = "'" + '0.00000000117'
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(49, 5).Value = '  -0.65%  '

# Row 50
$ws.Cells.Item(50, 4).Value = "'" + '1.701'
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 5).Value = '  +1.36%  '

# Row 51
$ws.Cells.Item(51, 4).Value = "'" + '9.008'
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Cells.Item(51, 5).Value = '  +0.02%  '
